$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the evaluation scores for row 14 (RES. 10)
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5
$ws.Range("F14").Value = 15
$ws.Range("G14").Value = 5
$ws.Range("H14").Value = 5
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 5
$ws.Range("K14").Value = 20
$ws.Range("L14").Value = 5
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 15
$ws.Range("P14").Value = 5
$ws.Range("Q14").Value = 5
$ws.Range("R14").Value = 5
$ws.Range("S14").Value = 5
$ws.Range("T14").Value = 20
$ws.Range("U14").Value = 5
$ws.Range("V14").Value = 5
$ws.Range("W14").Value = 5
$ws.Range("X14").Value = 15
$ws.Range("Y14").Value = 5
$ws.Range("Z14").Value = 5
$ws.Range("AA14").Value = 5
$ws.Range("AB14").Value = 5
$ws.Range("AC14").Value = 20
$ws.Range("AD14").Value = 5
$ws.Range("AE14").Value = 5

# Update the selected cell to match the saved view state
$ws.Range("O36").Select()
